$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93, shifting rows 93:145 down to 94:146
# (mirrors a new daily price record being inserted into the weekly log).
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new record's data.
$ws.Range("A93").Value = 7
$ws.Range("B93").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C93").Value = "Ñuble"
$ws.Range("D93").Value = 44518
$ws.Range("E93").Value = 16
$ws.Range("F93").Value = 100112032
$ws.Range("G93").Value = "Zapallo italiano"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 100
$ws.Range("K93").Value = 8000
$ws.Range("L93").Value = 9000
$ws.Range("M93").Value = 8500
$ws.Range("N93").Value = "$/caja 60 unidades"
$ws.Range("O93").Value = "Región del Maule"
$ws.Range("P93").Value = 142
$ws.Range("Q93").Value = 60
$ws.Range("R93").Value = "Hortaliza"
